$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: Install date for every collar (all rows share the same date
# and the same date style) ---
$ws.Range("E2").Value = 43217
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy()
$ws.Range("E3:E26").PasteSpecial(-4122)
$ws.Range("E3:E26").Value = 43217
$excel.CutCopyMode = 0
$ws.Columns("E:E").AutoFit() | Out-Null

# --- Column headers (order matters for shared-string indices: "Duration" must
# be registered before "Install date") ---
$ws.Range("F1").Value = "Duration"
$ws.Range("E1").Value = "Install date"

# --- Column F: Duration = retrieval date - install date ---
$ws.Range("F2").Formula = "=D2-E2"
$ws.Range("F3:F26").Formula = "=D3-E3"
$ws.Range("F2:F26").ClearFormats()
# Row 16's retrieval date is text ("XXXX"), so the duration there is a plain
# literal instead of a formula.
$ws.Range("F16").Value = 0

# --- Column G: secondary copy of the duration for a subset of rows ---
$ws.Range("G2").Formula = "=F2"
$ws.Range("G3:G26").Formula = "=F3"
$ws.Range("G2:G26").ClearFormats()
$ws.Range("G5").Clear()
$ws.Range("G6").Clear()
$ws.Range("G8").Clear()
$ws.Range("G10").Clear()
$ws.Range("G16").Clear()
$ws.Range("G17").Clear()
$ws.Range("G22").Clear()
$ws.Range("G23").Clear()
$ws.Range("G26").Clear()

# Match the selection left behind in the saved workbook
$ws.Range("G2:G26").Select()
